# structure revamp for future tools and internationalization
# Lowercase the "document" (h1/title) column values in columns A and F,
# rows 2-21, leaving the slug (column B) and Title-Cased (column C) columns
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 21; $row++) {
    $cellA = $ws.Cells.Item($row, 1)
    if ($cellA.Value2 -ne $null) {
        $cellA.Value2 = $cellA.Value2.ToString().ToLower()
    }

    $cellF = $ws.Cells.Item($row, 6)
    if ($cellF.Value2 -ne $null) {
        $cellF.Value2 = $cellF.Value2.ToString().ToLower()
    }
}
